$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"
$ws.Range("A2").Value = "akshata"
$ws.Range("B2").Value = "akshata"
$ws.Range("A3").Value = "adarsh"
$ws.Range("B3").Value = "adarsh"

$ws.Range("A5").Select()
